# Atualização de bases das ligas, do dia: 23-02-2024 às 08:18
#
# The upstream feed re-matched a handful of rows to the correct fixture
# (teams/odds/result columns B and F..AC swapped between rows while the
# row-id column A and the league/date columns C,D,E stay put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-match payload (everything except the row
# index in A and the constant Div/Div Original Name/Date columns C/D/E).
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Get-RowValues($row) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range($c + $row).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($c in $cols) {
        $ws.Range($c + $row).Value2 = $vals[$c]
    }
}

function Swap-Rows($rowA, $rowB) {
    $a = Get-RowValues $rowA
    $b = Get-RowValues $rowB
    Set-RowValues $rowA $b
    Set-RowValues $rowB $a
}

function Rotate-Rows($rows) {
    # rows[0] <- rows[1] <- rows[2] <- ... <- rows[0]
    $first = Get-RowValues $rows[0]
    $lastIndex = $rows.Length - 1
    for ($i = 0; $i -lt $lastIndex; $i++) {
        $nextRow = $rows[$i + 1]
        $nextVals = Get-RowValues $nextRow
        $curRow = $rows[$i]
        Set-RowValues $curRow $nextVals
    }
    $tailRow = $rows[$lastIndex]
    Set-RowValues $tailRow $first
}

# Row 231 <-> Row 232 (ids 229/230) were swapped.
Swap-Rows 231 232

# Rows 238, 239, 241 (ids 236/237/239) were rotated: 238<-241, 241<-239, 239<-238.
Rotate-Rows @(238, 241, 239)
